# Update NATMI LR-pair TPM output (Col8a1-Itga2) with new TPM-derived statistics.
# Ligand-expressing / receptor-expressing cell counts and all downstream
# expression / specificity / edge-weight columns (E:T) are refreshed for rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 3.0
$ws.Cells.Item(2, 6).Value2 = 1.0
$ws.Cells.Item(2, 7).Value2 = 2.502470333333333
$ws.Cells.Item(2, 8).Value2 = 7.507410999999999
$ws.Cells.Item(2, 9).Value2 = 0.07821409705091072
$ws.Cells.Item(2, 10).Value2 = 0.07821409705091072
$ws.Cells.Item(2, 11).Value2 = 3.0
$ws.Cells.Item(2, 12).Value2 = 1.0
$ws.Cells.Item(2, 13).Value2 = 0.229822
$ws.Cells.Item(2, 14).Value2 = 0.689466
$ws.Cells.Item(2, 15).Value2 = 0.09226175421862418
$ws.Cells.Item(2, 16).Value2 = 0.09226175421862419
$ws.Cells.Item(2, 17).Value2 = 0.5751227369473333
$ws.Cells.Item(2, 18).Value2 = 5.176104632526
$ws.Cells.Item(2, 19).Value2 = 0.007216169798542742
$ws.Cells.Item(2, 20).Value2 = 0.007216169798542744

$ws.Cells.Item(3, 5).Value2 = 3.0
$ws.Cells.Item(3, 6).Value2 = 1.0
$ws.Cells.Item(3, 7).Value2 = 2.502470333333333
$ws.Cells.Item(3, 8).Value2 = 7.507410999999999
$ws.Cells.Item(3, 9).Value2 = 0.07821409705091072
$ws.Cells.Item(3, 10).Value2 = 0.07821409705091072
$ws.Cells.Item(3, 15).Value2 = 0.4364142651333466
$ws.Cells.Item(3, 16).Value2 = 0.4364142651333466
$ws.Cells.Item(3, 17).Value2 = 2.720431328582777
$ws.Cells.Item(3, 18).Value2 = 24.483881957245
$ws.Cells.Item(3, 19).Value2 = 0.03413374768754145
$ws.Cells.Item(3, 20).Value2 = 0.03413374768754145

$ws.Cells.Item(4, 5).Value2 = 3.0
$ws.Cells.Item(4, 6).Value2 = 1.0
$ws.Cells.Item(4, 7).Value2 = 2.502470333333333
$ws.Cells.Item(4, 8).Value2 = 7.507410999999999
$ws.Cells.Item(4, 9).Value2 = 0.07821409705091072
$ws.Cells.Item(4, 10).Value2 = 0.07821409705091072
$ws.Cells.Item(4, 13).Value2 = 1.174057666666666
$ws.Cells.Item(4, 15).Value2 = 0.4713239806480292
$ws.Cells.Item(4, 16).Value2 = 0.4713239806480293
$ws.Cells.Item(4, 17).Value2 = 2.938044480455888
$ws.Cells.Item(4, 18).Value2 = 26.44240032410299
$ws.Cells.Item(4, 19).Value2 = 0.03686417956482652
$ws.Cells.Item(4, 20).Value2 = 0.03686417956482653

$ws.Cells.Item(5, 9).Value2 = 0.8193892102022395
$ws.Cells.Item(5, 10).Value2 = 0.8193892102022395
$ws.Cells.Item(5, 11).Value2 = 3.0
$ws.Cells.Item(5, 12).Value2 = 1.0
$ws.Cells.Item(5, 13).Value2 = 0.229822
$ws.Cells.Item(5, 14).Value2 = 0.689466
$ws.Cells.Item(5, 15).Value2 = 0.09226175421862418
$ws.Cells.Item(5, 16).Value2 = 0.09226175421862419
$ws.Cells.Item(5, 17).Value2 = 6.025120572444666
$ws.Cells.Item(5, 18).Value2 = 54.226085152002
$ws.Cells.Item(5, 19).Value2 = 0.07559828592107161
$ws.Cells.Item(5, 20).Value2 = 0.07559828592107162

$ws.Cells.Item(6, 9).Value2 = 0.8193892102022395
$ws.Cells.Item(6, 10).Value2 = 0.8193892102022395
$ws.Cells.Item(6, 15).Value2 = 0.4364142651333466
$ws.Cells.Item(6, 16).Value2 = 0.4364142651333466
$ws.Cells.Item(6, 17).Value2 = 28.49987613212388
$ws.Cells.Item(6, 19).Value2 = 0.3575931400286036
$ws.Cells.Item(6, 20).Value2 = 0.3575931400286036

$ws.Cells.Item(7, 9).Value2 = 0.8193892102022395
$ws.Cells.Item(7, 10).Value2 = 0.8193892102022395
$ws.Cells.Item(7, 13).Value2 = 1.174057666666666
$ws.Cells.Item(7, 15).Value2 = 0.4713239806480292
$ws.Cells.Item(7, 16).Value2 = 0.4713239806480293
$ws.Cells.Item(7, 17).Value2 = 30.77964250885343
$ws.Cells.Item(7, 18).Value2 = 277.0167825796809
$ws.Cells.Item(7, 19).Value2 = 0.3861977842525643
$ws.Cells.Item(7, 20).Value2 = 0.3861977842525644

$ws.Cells.Item(8, 7).Value2 = 3.276195666666666
$ws.Cells.Item(8, 8).Value2 = 9.828586999999999
$ws.Cells.Item(8, 9).Value2 = 0.1023966927468496
$ws.Cells.Item(8, 10).Value2 = 0.1023966927468496
$ws.Cells.Item(8, 11).Value2 = 3.0
$ws.Cells.Item(8, 12).Value2 = 1.0
$ws.Cells.Item(8, 13).Value2 = 0.229822
$ws.Cells.Item(8, 14).Value2 = 0.689466
$ws.Cells.Item(8, 15).Value2 = 0.09226175421862418
$ws.Cells.Item(8, 16).Value2 = 0.09226175421862419
$ws.Cells.Item(8, 17).Value2 = 0.7529418405046666
$ws.Cells.Item(8, 18).Value2 = 6.776476564542
$ws.Cells.Item(8, 19).Value2 = 0.00944729849900982
$ws.Cells.Item(8, 20).Value2 = 0.009447298499009821

$ws.Cells.Item(9, 7).Value2 = 3.276195666666666
$ws.Cells.Item(9, 8).Value2 = 9.828586999999999
$ws.Cells.Item(9, 9).Value2 = 0.1023966927468496
$ws.Cells.Item(9, 10).Value2 = 0.1023966927468496
$ws.Cells.Item(9, 15).Value2 = 0.4364142651333466
$ws.Cells.Item(9, 16).Value2 = 0.4364142651333466
$ws.Cells.Item(9, 17).Value2 = 3.561546848907222
$ws.Cells.Item(9, 18).Value2 = 32.05392164016499
$ws.Cells.Item(9, 19).Value2 = 0.04468737741720147
$ws.Cells.Item(9, 20).Value2 = 0.04468737741720147

$ws.Cells.Item(10, 7).Value2 = 3.276195666666666
$ws.Cells.Item(10, 8).Value2 = 9.828586999999999
$ws.Cells.Item(10, 9).Value2 = 0.1023966927468496
$ws.Cells.Item(10, 10).Value2 = 0.1023966927468496
$ws.Cells.Item(10, 13).Value2 = 1.174057666666666
$ws.Cells.Item(10, 15).Value2 = 0.4713239806480292
$ws.Cells.Item(10, 16).Value2 = 0.4713239806480293
$ws.Cells.Item(10, 17).Value2 = 2.938044480455888
$ws.Cells.Item(10, 18).Value2 = 26.44240032410299
$ws.Cells.Item(10, 19).Value2 = 0.03686417956482652
$ws.Cells.Item(10, 20).Value2 = 0.03686417956482653
